$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 96
$ws1.Range("F10").Value = 3517
$ws1.Range("F14").Value = 2805
$ws1.Range("F16").Value = 509
$ws1.Range("F17").Value = 2171
$ws1.Range("F23").Value = 142
$ws1.Range("F34").Value = 3580
$ws1.Range("F35").Value = 3123
$ws1.Range("F41").Value = 1293
$ws1.Range("F47").Value = 9

# Sheet 2: 演出 (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 179
$ws2.Range("F16").Value = 178

# Sheet 3: 本地生活 (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 745
$ws3.Range("F4").Value = 125
$ws3.Range("F5").Value = 2001

# Sheet 4: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 125
$ws4.Range("F9").Value = 96
$ws4.Range("F16").Value = 509
$ws4.Range("F17").Value = 2171
$ws4.Range("F21").Value = 142
$ws4.Range("F33").Value = 3581
$ws4.Range("F34").Value = 3123
$ws4.Range("F48").Value = 9
$ws4.Range("F49").Value = 178
